$d = $word.ActiveDocument

$d.Content.Find.Execute("83÷9=9, 2", $true, $false, $false, $false, $false, $true, 1, $false, "49÷7=7, 0", 2) | Out-Null
$d.Content.Find.Execute("24÷2=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "92÷4=23, 0", 2) | Out-Null
$d.Content.Find.Execute("67÷9=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "84÷9=9, 3", 2) | Out-Null
$d.Content.Find.Execute("21÷7=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "12÷9=1, 3", 2) | Out-Null
$d.Content.Find.Execute("42÷6=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "10÷5=2, 0", 2) | Out-Null
$d.Content.Find.Execute("54÷5=10, 4", $true, $false, $false, $false, $false, $true, 1, $false, "30÷5=6, 0", 2) | Out-Null
$d.Content.Find.Execute("52÷2=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "73÷8=9, 1", 2) | Out-Null
$d.Content.Find.Execute("87÷5=17, 2", $true, $false, $false, $false, $false, $true, 1, $false, "68÷8=8, 4", 2) | Out-Null
$d.Content.Find.Execute("72÷4=18, 0", $true, $false, $false, $false, $false, $true, 1, $false, "14÷9=1, 5", 2) | Out-Null
$d.Content.Find.Execute("44÷4=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "55÷7=7, 6", 2) | Out-Null
$d.Content.Find.Execute("20÷3=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "70÷2=35, 0", 2) | Out-Null
$d.Content.Find.Execute("36÷7=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "69÷4=17, 1", 2) | Out-Null
$d.Content.Find.Execute("32÷8=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "57÷3=19, 0", 2) | Out-Null
$d.Content.Find.Execute("49÷6=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "36÷6=6, 0", 2) | Out-Null
$d.Content.Find.Execute("31÷7=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "92÷7=13, 1", 2) | Out-Null
$d.Content.Find.Execute("65÷3=21, 2", $true, $false, $false, $false, $false, $true, 1, $false, "85÷6=14, 1", 2) | Out-Null
$d.Content.Find.Execute("18÷8=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "78÷8=9, 6", 2) | Out-Null
$d.Content.Find.Execute("21÷3=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "62÷2=31, 0", 2) | Out-Null
$d.Content.Find.Execute("36÷3=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "18÷3=6, 0", 2) | Out-Null
$d.Content.Find.Execute("17÷9=1, 8", $true, $false, $false, $false, $false, $true, 1, $false, "76÷5=15, 1", 2) | Out-Null
$d.Content.Find.Execute("52÷8=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "50÷4=12, 2", 2) | Out-Null
$d.Content.Find.Execute("82÷3=27, 1", $true, $false, $false, $false, $false, $true, 1, $false, "88÷5=17, 3", 2) | Out-Null
$d.Content.Find.Execute("47÷5=9, 2", $true, $false, $false, $false, $false, $true, 1, $false, "20÷2=10, 0", 2) | Out-Null
$d.Content.Find.Execute("11÷3=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "18÷3=6, 0", 2) | Out-Null
$d.Content.Find.Execute("38÷5=7, 3", $true, $false, $false, $false, $false, $true, 1, $false, "69÷3=23, 0", 2) | Out-Null
